# Append daily COVID-style figures for Modena through 2021-06-28 (serial 44375),
# continuing the existing A:D table that previously ended at row 269 (serial 44343).
# Columns: A = date serial, B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(44344, 15, 110, 58.19705522900541),
    @(44345, 18, 100, 52.90641384455037),
    @(44346, 11, 111, 58.72611936745092),
    @(44347, 10, 88, 46.55764418320432),
    @(44348, 2, 83, 43.91232349097681),
    @(44349, 4, 81, 42.8541952140858),
    @(44350, 0, 60, 31.74384830673023),
    @(44351, 9, 54, 28.5694634760572),
    @(44352, 0, 36, 19.04630898403813),
    @(44353, 7, 32, 16.93005243025612),
    @(44354, 13, 35, 18.51724484559263),
    @(44355, 3, 36, 19.04630898403813),
    @(44356, 3, 35, 18.51724484559263),
    @(44357, 5, 40, 21.16256553782015),
    @(44358, 4, 35, 18.51724484559263),
    @(44359, 4, 39, 20.63350139937465),
    @(44360, 2, 34, 17.98818070714713),
    @(44361, 11, 32, 16.93005243025612),
    @(44362, 0, 29, 15.34286001491961),
    @(44363, 1, 27, 14.2847317380286),
    @(44364, 1, 23, 12.16847518424659),
    @(44365, 0, 19, 10.05221863046457),
    @(44366, 8, 23, 12.16847518424659),
    @(44367, 1, 22, 11.63941104580108),
    @(44368, 4, 15, 7.935962076682557),
    @(44369, 4, 19, 10.05221863046457),
    @(44370, 3, 21, 11.11034690735558),
    @(44371, 2, 22, 11.63941104580108),
    @(44372, 4, 26, 13.7556675995831),
    @(44373, 10, 28, 14.81379587647411),
    @(44374, 5, 32, 16.93005243025612),
    @(44375, 4, 32, 16.93005243025612)
)

$firstNewRow = 270

# The date column (A) uses style index 2 (center/top aligned, bordered, custom date
# number format). Grab that formatting from the last existing row so the appended
# cells match, then paste values on top row by row.
$styleSource = $ws.Cells.Item($firstNewRow - 1, 1)
$styleSource.Copy()

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $values = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $values[0]
    $dateCell.PasteSpecial(-4122)  # xlPasteFormats - copy date style, keep the value just set

    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
}

$excel.CutCopyMode = 0

Write-Host "Wrote rows $firstNewRow to $($firstNewRow + $newRows.Count - 1)"
